$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 217 (pushes existing rows 217-299 down to 218-300)
$ws.Rows(217).Insert()

# Populate the newly inserted row with the new price-report entry
$ws.Range("A217").Value = 7
$ws.Range("B217").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C217").Value = "Ñuble"
$ws.Range("D217").Value = 44988
$ws.Range("E217").Value = 16
$ws.Range("F217").Value = 100112043
$ws.Range("G217").Value = "Pepino ensalada"
$ws.Range("H217").Value = "Sin especificar"
$ws.Range("I217").Value = "Primera"
$ws.Range("J217").Value = 120
$ws.Range("K217").Value = 9000
$ws.Range("L217").Value = 9000
$ws.Range("M217").Value = 9000
$ws.Range("N217").Value = "$/caja 80 unidades"
$ws.Range("O217").Value = "Región del Maule"
$ws.Range("P217").Value = 112
$ws.Range("Q217").Value = 80
$ws.Range("R217").Value = "Hortaliza"
